# Update gh-pages output: refresh exhibition counters and insert a new
# exhibition entry ("常熟·CDW·动漫展03") on 2024-05-03 into both the
# "展览" (sheet 1) and "全部类型" (sheet 4) worksheets.

$wb = $excel.ActiveWorkbook

foreach ($sheetIndex in 1, 4) {
    $ws = $wb.Worksheets.Item($sheetIndex)

    # --- simple scalar "want-to-go" counter refreshes -------------------
    $ws.Range("F4").Value  = 1595
    $ws.Range("F5").Value  = 612
    $ws.Range("F8").Value  = 11403
    $ws.Range("F9").Value  = 21
    $ws.Range("F12").Value = 350
    $ws.Range("F14").Value = 789
    $ws.Range("F15").Value = 12350
    $ws.Range("F16").Value = 13015

    # --- insert a brand-new row 21, pushing the old rows 21-23 down to ---
    # --- 22-24 -------------------------------------------------------------
    $ws.Rows.Item(21).Insert()

    # Re-apply the header/index-column formatting (bold, centered, boxed)
    # that is normally carried by column A so the new row matches its
    # neighbours instead of whatever Excel guessed during the insert.
    $ws.Range("A20").Copy() | Out-Null
    $ws.Range("A21").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
    $excel.CutCopyMode = $false

    # The "start date" column stores plain text such as "2024-05-03"
    # rather than a real date value, so prefix it so Excel keeps it as
    # literal text instead of reinterpreting it as a date serial, then
    # drop the formatting quirk that introduces so the cell stays
    # styleless like its neighbours.
    $ws.Range("A21").Value = 20
    $ws.Range("B21").Value = "'2024-05-03"
    $ws.Range("B21").ClearFormats()
    $ws.Range("C21").Value = "常熟·CDW·动漫展03"
    $ws.Range("D21").Value = "常熟国际展览中心 国际展览中心"
    $ws.Range("E21").Value = "2024.05.03 09:00-05.04 17:30"
    $ws.Range("F21").Value = 6
    $ws.Range("G21").Value = 60
    $ws.Range("H21").Value = "https://show.bilibili.com/platform/detail.html?id=82489"
    $ws.Range("I21").Value = "//i0.hdslb.com/bfs/openplatform/202403/XK411blC1709794808211.jpeg"

    # The row-index column (A) numbers every record sequentially and is
    # independent of the data shift above, so renumber rows 22-24 back to
    # 21, 22, 23 after the insert shoved the old rows' index values down
    # along with everything else.
    $ws.Range("A22").Value = 21
    $ws.Range("A23").Value = 22
    $ws.Range("A24").Value = 23
}
